# Generate Report for Archive
#
# Localization status has moved on: items that were previously
# "Ready for handoff" are now "In Translation". Update the status text
# everywhere it appears, then tighten the now-narrower status columns
# to match.

$wb = $excel.ActiveWorkbook

# --- 1. Update status text on every sheet (Overview, zh-cn, de-de) ---
foreach ($ws in $wb.Worksheets) {
    $ws.Cells.Replace("Ready for handoff", "In Translation") | Out-Null
}

# --- 2. Re-fit the status columns now that the text is shorter ---
# Overview sheet: status columns are E (zh-cn) and F (de-de)
$overview = $wb.Worksheets.Item("Overview")
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

# zh-cn sheet: status column is C
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Columns.Item(3).ColumnWidth = 12.5

# de-de sheet: status column is C
$dede = $wb.Worksheets.Item("de-de")
$dede.Columns.Item(3).ColumnWidth = 12.5
